# Update Upload routine - overwrite cover sheet details and lane details
$wb = $excel.ActiveWorkbook

# Cover sheet is the 1st worksheet tab
$coverSheet = $wb.Worksheets.Item(1)
$coverSheet.Range("A2").Value = "ABC Bowling Center"
$coverSheet.Range("A3").Value = "Arlington, TX"

# Lane details sheet is the 7th worksheet tab
$laneSheet = $wb.Worksheets.Item(7)
$laneSheet.Range("A4").Value = "Lane Suface: Wood"
$laneSheet.Range("A5").Value = "Year  of Installation: 2000"
$laneSheet.Range("A7").Value = "Underlayment Year Installation: 2001"
$laneSheet.Range("A10").Value = "Pin Decks: DBA"
